$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Extend the table (Table2) by one row so its range grows from A2:G28 to
# A2:G29 (matching xl/tables/table1.xml + autoFilter change in the diff).
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# ---------------------------------------------------------------------------
# Pre-format the new rows (alignment only, before any values are written) so
# that the resulting style table reuses the existing entries wherever
# possible and the shared-string table fills up in authoring order.
# ---------------------------------------------------------------------------
$ws.Range("A9:E9").HorizontalAlignment = $xlCenter
$ws.Range("A9:E9").VerticalAlignment = $xlCenter
$ws.Range("F9:G9").HorizontalAlignment = $xlCenter
$ws.Range("F9:G9").VerticalAlignment = $xlCenter

$ws.Range("A10:E10").HorizontalAlignment = $xlCenter
$ws.Range("A10:E10").VerticalAlignment = $xlCenter
$ws.Range("F10").HorizontalAlignment = $xlCenter
$ws.Range("F10").VerticalAlignment = $xlCenter
$ws.Range("G10").HorizontalAlignment = $xlCenter
$ws.Range("G10").VerticalAlignment = $xlCenter
$ws.Range("G10").WrapText = $true

$ws.Range("A11").HorizontalAlignment = $xlCenter
$ws.Range("B11").HorizontalAlignment = $xlCenter
$ws.Range("C11").HorizontalAlignment = $xlCenter
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("E11").HorizontalAlignment = $xlCenter
$ws.Range("F11:G11").HorizontalAlignment = $xlCenter
$ws.Range("F11:G11").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 10 : SendtoDaq() / {1 ,8 ,7 ,5, 4 , 5, 9,14}(uint16_t) / 8 / ...
# (entered first, matching the workbook's authoring order)
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "{1 ,8 ,7 ,5, 4 , 5, 9,14}(uint16_t)"
$ws.Range("A10").Value = "SendtoDaq()"
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = "Putty Output: position to go to decrementing"
$ws.Range("E10").Value = "Putty Output: position to go to decrementing"
$ws.Range("F10").Value = "'---"
$ws.Range("G10").Value = "Because we are giving it a constant value" + [char]10 + "the error will keep on decrementing " + [char]10 + "to control the motor" + [char]10

$ws.Rows.Item(10).RowHeight = 90

# ---------------------------------------------------------------------------
# Row 9 : getPosition() / {1 ,8 ,7 ,5, 4 , 5, 9,14} (uint16_t) / 7 / ...
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "getPosition()"
$ws.Range("B9").Value = "{1 ,8 ,7 ,5, 4 , 5, 9,14} (uint16_t)"
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = "Putty Output: position to go to is: 33000"
$ws.Range("E9").Value = "Putty Output: position to go to is: 33000"
$ws.Range("F9").Value = "'---"
$ws.Range("G9").Value = "'---"

# ---------------------------------------------------------------------------
# Row 11 : Gather_Data() / --- / 9 / Putty Output Showing a binary count...
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Gather_Data()"
$ws.Range("B11").Value = "'---"
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = "Putty Output Showing a binary count going from 000 --> 111 and then reset "
$ws.Range("E11").Value = "Putty Output Showing a binary count going from 000 --> 111 and then reset "
$ws.Range("F11").Value = "'---"
$ws.Range("G11").Value = "'---"

# ---------------------------------------------------------------------------
# Column D:E width grows from 39.71 to 68.71 characters (closest achievable
# width is used because of the host's width quantization).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 67.8
$ws.Columns.Item(5).ColumnWidth = 67.8

# ---------------------------------------------------------------------------
# Update the view: scroll back to the top and select D10 (matching the
# sheetView/selection change in the diff).
# ---------------------------------------------------------------------------
$ws.Range("D10").Select() | Out-Null

$wb.Save()
